# Update build/version timestamp strings across the workbook.
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

# --- "About" sheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: " + $newVersion

$wsAbout.Range("A6").Value = "Recommended Citation:  " + '"Global Energy Monitor, Coal mine boundaries and methane sources for Blackwater Coal Mine, Australia, M0011, version ' + "'" + $newVersion + "'" + ". (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 23; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S
    if ($cell.Value() -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
